$wb = $excel.ActiveWorkbook

# Fill in the research flag ("" ) for rows 6-14 in column Q of the Data sheet
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Range("Q6:Q14").Value = '"" '

# Move the active tab / selection from Research to Data, like a user who
# clicked over to the Data sheet and selected cell R5
$dataSheet.Activate() | Out-Null
$dataSheet.Range("R5").Select() | Out-Null
